# fixed bug change alias key config program error.
# - Update the "Sheet1 (3)" sheet: column A keys C1/C2/C3/A1/A2/A3 -> M1/M2/M3/N1/N2/N3
# - Update the "Sheet1 (2)" sheet: column A keys C1/C2/C3 -> M1/M2/M3 (A5:A7 left untouched)
# - Reorder sheet tabs so "Sheet1 (3)" comes before "Sheet1 (2)"

$wb = $excel.ActiveWorkbook

# --- Edit "Sheet1 (3)" first (alias keys + child keys) ---
$ws3 = $wb.Worksheets.Item("Sheet1 (3)")
$ws3.Activate()
$ws3.Range("A2").Value = "M1"
$ws3.Range("A3").Value = "M2"
$ws3.Range("A4").Value = "M3"
$ws3.Range("A5").Value = "N1"
$ws3.Range("A6").Value = "N2"
$ws3.Range("A7").Value = "N3"
$ws3.Range("C9").Select()

# --- Edit "Sheet1 (2)" (alias keys only) ---
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Activate()
$ws2.Range("A2").Value = "M1"
$ws2.Range("A3").Value = "M2"
$ws2.Range("A4").Value = "M3"
$ws2.Range("A2:A4").Select()

# --- Reorder tabs: move "Sheet1 (3)" before "Sheet1 (2)" ---
$ws3.Move($ws2)
